$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was updated from
# 2023-10-04 (45203) to 2023-10-06 (45205) for every data row (rows 2-92).
for ($row = 2; $row -le 92; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value = 45205
    }
}
